# Modeling and Prediction.pptx -- apply the committed edit:
#  1. Remove the trailing blank slide (slide id 265 / last slide in the deck).
#  2. Refresh the cached "datetimeFigureOut" footer-date text (master + all
#     layouts) from 3/24/2020 -> 3/31/2020, matching PowerPoint's automatic
#     field re-cache on save.

$p = $ppt.ActivePresentation

# --- 1. Delete the trailing empty slide -------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Delete()

# --- 2. Update the cached date-field text ------------------------------
$newDate = "3/31/2020"

function Update-DatePlaceholder($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.Name -like "Date Placeholder*") {
      if ($sh.TextFrame.TextRange.Text -ne $newDate) {
        $sh.TextFrame.TextRange.Text = $newDate
      }
    }
  }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
  $layout = $layouts.Item($li)
  Update-DatePlaceholder $layout.Shapes
}
